$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Cells.Item(94, 1) 'Páprica Defumada Kitano Reserva 36g'
$ws.Cells.Item(94, 2).Value = 0
Set-TextCell $ws.Cells.Item(94, 3) 'https://mercado.carrefour.com.br/farofa-de-mandioca-tradicional-yoki-400g-6582613/p'
Set-TextCell $ws.Cells.Item(94, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(95, 1) 'Não encontrado'
$ws.Cells.Item(95, 2).Value = 0
Set-TextCell $ws.Cells.Item(95, 3) 'https://mercado.carrefour.com.br/massa-para-pastel-discao-massa-leve-500g-841757/p'
Set-TextCell $ws.Cells.Item(95, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(96, 1) 'Pimentão Block Vermelho Trebeshi 150 g'
$ws.Cells.Item(96, 2).Value = 0
Set-TextCell $ws.Cells.Item(96, 3) 'https://mercado.carrefour.com.br/pimentao-block-vermelho-trebeshi-150-g-5738458/p'
Set-TextCell $ws.Cells.Item(96, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(97, 1) 'Não encontrado'
$ws.Cells.Item(97, 2).Value = 0
Set-TextCell $ws.Cells.Item(97, 3) 'https://mercado.carrefour.com.br/chocolate-ao-leite-com-amendoim-shot-165g-5790859/p'
Set-TextCell $ws.Cells.Item(97, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(98, 1) 'Couve Flor Cledson 300 g'
$ws.Cells.Item(98, 2).Value = 0
Set-TextCell $ws.Cells.Item(98, 3) 'https://mercado.carrefour.com.br/couve-flor-cledson-300-g-9560297/p'
Set-TextCell $ws.Cells.Item(98, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(99, 1) 'Não encontrado'
$ws.Cells.Item(99, 2).Value = 0
Set-TextCell $ws.Cells.Item(99, 3) 'https://mercado.carrefour.com.br/patinho-fracionado-a-vacuo-500g-18325/p'
Set-TextCell $ws.Cells.Item(99, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(100, 1) 'Não encontrado'
$ws.Cells.Item(100, 2).Value = 0
Set-TextCell $ws.Cells.Item(100, 3) 'https://mercado.carrefour.com.br/paleta-bovina-a-vacuo-500gnao-reativarcodigo-de-compra-20745/p'
Set-TextCell $ws.Cells.Item(100, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(101, 1) 'Não encontrado'
$ws.Cells.Item(101, 2).Value = 0
Set-TextCell $ws.Cells.Item(101, 3) 'https://mercado.carrefour.com.br/costela-minga-bovina-cong-aprox-2kg-224006/p'
Set-TextCell $ws.Cells.Item(101, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(102, 1) 'Não encontrado'
$ws.Cells.Item(102, 2).Value = 0
Set-TextCell $ws.Cells.Item(102, 3) 'https://mercado.carrefour.com.br/presunto-cozido-sem-capa-fatiado-aurora-aproximadamente-200-g-49450/p'
Set-TextCell $ws.Cells.Item(102, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(103, 1) 'Não encontrado'
$ws.Cells.Item(103, 2).Value = 0
Set-TextCell $ws.Cells.Item(103, 3) 'https://mercado.carrefour.com.br/mortadela-defumada-sadia-280g-5447045/p'
Set-TextCell $ws.Cells.Item(103, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(104, 1) 'Não encontrado'
$ws.Cells.Item(104, 2).Value = 0
Set-TextCell $ws.Cells.Item(104, 3) 'https://mercado.carrefour.com.br/queijo-minas-frescal-aurora-450-g-6264693/p'
Set-TextCell $ws.Cells.Item(104, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(105, 1) 'Não encontrado'
$ws.Cells.Item(105, 2).Value = 0
Set-TextCell $ws.Cells.Item(105, 3) 'https://mercado.carrefour.com.br/queijo-coalho-bom-leite-500-g-4305054/p'
Set-TextCell $ws.Cells.Item(105, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(106, 1) 'Não encontrado'
$ws.Cells.Item(106, 2).Value = 0
Set-TextCell $ws.Cells.Item(106, 3) 'https://mercado.carrefour.com.br/leite-uht-integral-piratininga-1-l-665017/p'
Set-TextCell $ws.Cells.Item(106, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(107, 1) 'Não encontrado'
$ws.Cells.Item(107, 2).Value = 0
Set-TextCell $ws.Cells.Item(107, 3) 'https://mercado.carrefour.com.br/iogurte-natural-tradicional-batavo-170g-5150439/p'
Set-TextCell $ws.Cells.Item(107, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(108, 1) 'Não encontrado'
$ws.Cells.Item(108, 2).Value = 0
Set-TextCell $ws.Cells.Item(108, 3) 'https://mercado.carrefour.com.br/manteiga-com-sal-aviacao-200-g-10010/p'
Set-TextCell $ws.Cells.Item(108, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(109, 1) 'Não encontrado'
$ws.Cells.Item(109, 2).Value = 0
Set-TextCell $ws.Cells.Item(109, 3) 'https://mercado.carrefour.com.br/creme-de-leite-ultrapasteurizado-itambe-200-g-5988921/p'
Set-TextCell $ws.Cells.Item(109, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(110, 1) 'Não encontrado'
$ws.Cells.Item(110, 2).Value = 0
Set-TextCell $ws.Cells.Item(110, 3) 'https://mercado.carrefour.com.br/requeijao-cremoso-aviacao-tradicional-220-g-10000/p'
Set-TextCell $ws.Cells.Item(110, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(111, 1) 'Não encontrado'
$ws.Cells.Item(111, 2).Value = 0
Set-TextCell $ws.Cells.Item(111, 3) 'https://mercado.carrefour.com.br/acucar-cristal-carrefour-1kg-5147300/p'
Set-TextCell $ws.Cells.Item(111, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(112, 1) 'Não encontrado'
$ws.Cells.Item(112, 2).Value = 0
Set-TextCell $ws.Cells.Item(112, 3) 'https://mercado.carrefour.com.br/mel-com-cacau-e-avela-400-g-4510146/p'
Set-TextCell $ws.Cells.Item(112, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(113, 1) 'Não encontrado'
$ws.Cells.Item(113, 2).Value = 0
Set-TextCell $ws.Cells.Item(113, 3) 'https://mercado.carrefour.com.br/geleia-de-goiaba-selecoes-c-pedacos-260-g-1280815/p'
Set-TextCell $ws.Cells.Item(113, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(114, 1) 'Não encontrado'
$ws.Cells.Item(114, 2).Value = 0
Set-TextCell $ws.Cells.Item(114, 3) 'https://mercado.carrefour.com.br/suco-de-uva-integral-maric-1-l-3538256/p'
Set-TextCell $ws.Cells.Item(114, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(115, 1) 'Não encontrado'
$ws.Cells.Item(115, 2).Value = 0
Set-TextCell $ws.Cells.Item(115, 3) 'https://mercado.carrefour.com.br/vinho-tinto-fino-seco-cabernet-sauvignon-pergola-750ml-1521709/p'
Set-TextCell $ws.Cells.Item(115, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(116, 1) 'Não encontrado'
$ws.Cells.Item(116, 2).Value = 0
Set-TextCell $ws.Cells.Item(116, 3) 'https://mercado.carrefour.com.br/whisky-red-label-johnnie-walker-1-litro-2719/p'
Set-TextCell $ws.Cells.Item(116, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(117, 1) 'Não encontrado'
$ws.Cells.Item(117, 2).Value = 0
Set-TextCell $ws.Cells.Item(117, 3) 'https://mercado.carrefour.com.br/refrigerante-coca-cola-sabor-cola-1-5-l-11087/p'
Set-TextCell $ws.Cells.Item(117, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(118, 1) 'Não encontrado'
$ws.Cells.Item(118, 2).Value = 0
Set-TextCell $ws.Cells.Item(118, 3) 'https://mercado.carrefour.com.br/cafe-torrado-e-moido-extraforte-melitta-500g-271203/p'
Set-TextCell $ws.Cells.Item(118, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(119, 1) 'Não encontrado'
$ws.Cells.Item(119, 2).Value = 0
Set-TextCell $ws.Cells.Item(119, 3) 'https://mercado.carrefour.com.br/farinha-de-trigo-dona-benta-tradicional-1kg-196416/p'
Set-TextCell $ws.Cells.Item(119, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(120, 1) 'Não encontrado'
$ws.Cells.Item(120, 2).Value = 0
Set-TextCell $ws.Cells.Item(120, 3) 'https://mercado.carrefour.com.br/azeite-extravirgem-portugues-oliveira-da-serra-500-ml-4526108/p'
Set-TextCell $ws.Cells.Item(120, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(121, 1) 'Não encontrado'
$ws.Cells.Item(121, 2).Value = 0
Set-TextCell $ws.Cells.Item(121, 3) 'https://mercado.carrefour.com.br/oleo-de-soja-soya-900ml-482616/p'
Set-TextCell $ws.Cells.Item(121, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(122, 1) 'Não encontrado'
$ws.Cells.Item(122, 2).Value = 0
Set-TextCell $ws.Cells.Item(122, 3) 'https://mercado.carrefour.com.br/margarina-qualy-com-sal-250g-4815618/p'
Set-TextCell $ws.Cells.Item(122, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(123, 1) 'Não encontrado'
$ws.Cells.Item(123, 2).Value = 0
Set-TextCell $ws.Cells.Item(123, 3) 'https://mercado.carrefour.com.br/arroz-branco-longofino-tipo-1-tio-joao-1kg-115658/p'
Set-TextCell $ws.Cells.Item(123, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(124, 1) 'Não encontrado'
$ws.Cells.Item(124, 2).Value = 0
Set-TextCell $ws.Cells.Item(124, 3) 'https://mercado.carrefour.com.br/feijao-preto-tipo-1-kicaldo-1kg-466510/p'
Set-TextCell $ws.Cells.Item(124, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(125, 1) 'Arroz Branco Carrefour Classic Olimpíadas 1Kg'
$ws.Cells.Item(125, 2).Value = 0
Set-TextCell $ws.Cells.Item(125, 3) 'https://mercado.carrefour.com.br/arroz-branco-carrefour-classic-olimpiadas-1kg-3433455/p'
Set-TextCell $ws.Cells.Item(125, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(126, 1) 'Não encontrado'
$ws.Cells.Item(126, 2).Value = 0
Set-TextCell $ws.Cells.Item(126, 3) 'https://mercado.carrefour.com.br/busca/pao%20frances'
Set-TextCell $ws.Cells.Item(126, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(127, 1) 'Não encontrado'
$ws.Cells.Item(127, 2).Value = 0
Set-TextCell $ws.Cells.Item(127, 3) 'https://mercado.carrefour.com.br/busca/biscoito%20doce'
Set-TextCell $ws.Cells.Item(127, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(128, 1) 'Não encontrado'
$ws.Cells.Item(128, 2).Value = 0
Set-TextCell $ws.Cells.Item(128, 3) 'https://mercado.carrefour.com.br/frango-inteiro-temperado-seara-assa-facil-aprox-19kg-170739/p'
Set-TextCell $ws.Cells.Item(128, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(129, 1) 'Não encontrado'
$ws.Cells.Item(129, 2).Value = 0
Set-TextCell $ws.Cells.Item(129, 3) 'https://mercado.carrefour.com.br/busca/cafe%20moido'
Set-TextCell $ws.Cells.Item(129, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(130, 1) 'Não encontrado'
$ws.Cells.Item(130, 2).Value = 0
Set-TextCell $ws.Cells.Item(130, 3) 'https://mercado.carrefour.com.br/busca/costela?page=1'
Set-TextCell $ws.Cells.Item(130, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(131, 1) 'Não encontrado'
$ws.Cells.Item(131, 2).Value = 0
Set-TextCell $ws.Cells.Item(131, 3) 'https://mercado.carrefour.com.br/costela-de-cordeiro-a-vacuo-28738/p'
Set-TextCell $ws.Cells.Item(131, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(132, 1) 'Não encontrado'
$ws.Cells.Item(132, 2).Value = 0
Set-TextCell $ws.Cells.Item(132, 3) 'https://mercado.carrefour.com.br/busca/lingui%C3%A7a'
Set-TextCell $ws.Cells.Item(132, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(133, 1) 'Não encontrado'
$ws.Cells.Item(133, 2).Value = 0
Set-TextCell $ws.Cells.Item(133, 3) 'https://mercado.carrefour.com.br/busca/lingui%C3%A7a?page=3'
Set-TextCell $ws.Cells.Item(133, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(134, 1) 'Ovo Branco Grande Mantiqueira Happy Eggs com 20 Unidades'
$ws.Cells.Item(134, 2).Value = 0
Set-TextCell $ws.Cells.Item(134, 3) 'https://mercado.carrefour.com.br/ovo-branco-grande-mantiqueira-happy-eggs-com-20-unidades-6403565/p'
Set-TextCell $ws.Cells.Item(134, 4) '2025-09-16'

Set-TextCell $ws.Cells.Item(135, 1) 'Óleo de Soja Confiare 900ml'
$ws.Cells.Item(135, 2).Value = 0
Set-TextCell $ws.Cells.Item(135, 3) 'https://mercado.carrefour.com.br/oleo-de-soja-confiare-900ml-3731243/p'
Set-TextCell $ws.Cells.Item(135, 4) '2025-09-16'
